$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column A (C0) and column C (idle), rows 2-51 (B / N column is unchanged: 1..50)
$AValues = @(100,203,298,401,502,595,698,798,891,984,1077,1166,1268,1370,1454,1549,1642,1735,1829,1917,2006,2111,2188,2295,2387,2496,2582,2668,2763,2857,2942,3023,3115,3186,3270,3361,3451,3548,3640,3725,3820,3905,3992,4077,4161,4256,4348,4417,4498,4608)
$CValues = @(66.180000000000007,66.12,66.06,66.010000000000005,65.95,65.89,65.83,65.78,65.72,65.66,65.599999999999994,65.55,65.489999999999995,65.430000000000007,65.37,65.319999999999993,65.260000000000005,65.2,65.150000000000006,65.09,65.03,64.98,64.92,64.87,64.81,64.760000000000005,64.7,64.650000000000006,64.59,64.540000000000006,64.48,64.430000000000007,64.37,64.319999999999993,64.260000000000005,64.209999999999994,64.150000000000006,64.099999999999994,64.040000000000006,63.99,63.93,63.88,63.83,63.77,63.72,63.66,63.61,63.56,63.5,63.45)

for ($i = 0; $i -lt $AValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $AValues[$i]
    $ws.Cells.Item($row, 3).Value = $CValues[$i]
}

# Update the sheet view: move the selection to J23 (also resets scroll/topLeftCell to default)
$ws.Range("J23").Select()
